$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report date range) ---
$ws.Range("A8").Value = "Volume 33   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/5/2026  Through  1/11/2026"

# --- Helper functions for cells whose underlying type/style must flip
#     between a numeric style and a text-placeholder style ("0" / "***.*") ---
function Set-NumericFromTemplate($cellRef, $templateRef, $num) {
    $ws.Range($templateRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122) # xlPasteFormats
    $ws.Range($cellRef).Value = $num
}

function Set-PlaceholderFromTemplate($cellRef, $templateRef, $text) {
    $ws.Range($templateRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122) # xlPasteFormats
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($templateRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122) # xlPasteFormats (restore exact style index)
}

# --- Style-swap cells ---
Set-PlaceholderFromTemplate "F15" "C14" "0"
Set-NumericFromTemplate "L20" "H14" -62.5
Set-NumericFromTemplate "C23" "G14" 1
Set-NumericFromTemplate "I23" "G14" 1
Set-PlaceholderFromTemplate "F27" "C14" "0"
Set-NumericFromTemplate "J28" "G14" 5
Set-NumericFromTemplate "K28" "H14" -20
Set-PlaceholderFromTemplate "D29" "C14" "0"
Set-PlaceholderFromTemplate "E29" "E14" "***.*"
Set-PlaceholderFromTemplate "D30" "C14" "0"
Set-PlaceholderFromTemplate "E30" "E14" "***.*"

# --- Plain numeric value updates ---
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -30
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = -11.111111111111
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = -14.285714285714
$ws.Range("L16").Value = -47.826086956521
$ws.Range("M16").Value = 33.333333333333
$ws.Range("N16").Value = -76
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 46
$ws.Range("G17").Value = 58
$ws.Range("H17").Value = -20.689655172413
$ws.Range("I17").Value = 18
$ws.Range("J17").Value = 22
$ws.Range("K17").Value = -18.181818181818
$ws.Range("L17").Value = -18.181818181818
$ws.Range("M17").Value = 260
$ws.Range("N17").Value = -25
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 16.666666666666
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 13
$ws.Range("J18").Value = 10
$ws.Range("K18").Value = 30
$ws.Range("L18").Value = 8.333333333333
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = -62.857142857142
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -15.384615384615
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = 15.789473684210
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 19
$ws.Range("K19").Value = -26.315789473684
$ws.Range("L19").Value = 16.666666666666
$ws.Range("M19").Value = 75
$ws.Range("N19").Value = -30
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 3
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = -57.142857142857
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -86.363636363636
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 46
$ws.Range("E21").Value = -15.217391304347
$ws.Range("F21").Value = 154
$ws.Range("G21").Value = 169
$ws.Range("H21").Value = -8.875739644970
$ws.Range("I21").Value = 60
$ws.Range("J21").Value = 72
$ws.Range("K21").Value = -16.666666666666
$ws.Range("L21").Value = -22.077922077922
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = -61.538461538461
$ws.Range("F23").Value = 2
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 12
$ws.Range("F24").Value = 88
$ws.Range("G24").Value = 116
$ws.Range("H24").Value = -24.137931034482
$ws.Range("I24").Value = 39
$ws.Range("J24").Value = 36
$ws.Range("K24").Value = 8.333333333333
$ws.Range("L24").Value = 30
$ws.Range("M24").Value = 77.272727272727
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -58.823529411764
$ws.Range("I25").Value = 7
$ws.Range("J25").Value = 9
$ws.Range("K25").Value = -22.222222222222
$ws.Range("L25").Value = -22.222222222222
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 12.5
$ws.Range("F26").Value = 73
$ws.Range("G26").Value = 72
$ws.Range("H26").Value = 1.388888888888
$ws.Range("I26").Value = 28
$ws.Range("J26").Value = 22
$ws.Range("K26").Value = 27.272727272727
$ws.Range("L26").Value = -15.151515151515
$ws.Range("M26").Value = 47.368421052631
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -40
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -27.272727272727
$ws.Range("I28").Value = 4
$ws.Range("J43").Value = 275
$ws.Range("K43").Value = -36.194895591647
$ws.Range("L43").Value = -58.893871449925
$ws.Range("M43").Value = -79.055597867479
$ws.Range("N43").Value = -83.139178418148
$ws.Range("J44").Value = 628
$ws.Range("K44").Value = 122.695035460993
$ws.Range("L44").Value = 115.068493150685
$ws.Range("M44").Value = 35.637149028077
$ws.Range("N44").Value = 38.021978021978
$ws.Range("J46").Value = 2400
$ws.Range("K46").Value = 24.416796267496
$ws.Range("L46").Value = -0.538748445917
$ws.Range("M46").Value = -44.814899977006
$ws.Range("N46").Value = -51.338199513382

Write-Host "edit complete"
